$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (rows 139-219) down by one row (140-220),
# preserving all column values (A:R) exactly as Excel's own fill/insert would.
$srcRange = $ws.Range("A139:R219")
$vals = $srcRange.Value2
$dstRange = $ws.Range("A140:R220")
$dstRange.Value = $vals

# The newly created row 220 did not inherit the date number format that
# column D uses elsewhere in the table - copy it from the row above.
$ws.Range("D220").NumberFormat = $ws.Range("D219").NumberFormat

# Insert the new record at row 139 (pushes the rest of the table down,
# already handled above). Only the fields that actually change are set;
# the remaining columns keep the values copied from the old row 139 above.
$ws.Range("D139").Value = 44606
$ws.Range("J139").Value = 480
$ws.Range("K139").Value = 17500
$ws.Range("L139").Value = 18000
$ws.Range("M139").Value = 17750
$ws.Range("P139").Value = 1775
